$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Level1")

# Update the displayed/active selection to C2
$ws.Range("C2").Select()

# Update the COUNT formula range in C1
$ws.Range("C1").Formula = "=COUNT(A2:A99)"

# Swap the "6" marker between R6 and Q7 (move to R8 / Q9)
$ws.Range("R6").Value = -1
$ws.Range("Q7").Value = -1
$ws.Range("R8").Value = 6
$ws.Range("Q9").Value = 6

# Row 11: H11 loses its "6" marker, P11:S11 drop from 0 to -1
$ws.Range("H11").Value = -1
$ws.Range("P11:S11").Value = -1

# Row 12: H12 gains the "6" marker, P12:S12 drop from 5 to 0
$ws.Range("H12").Value = 6
$ws.Range("P12:S12").Value = 0
